$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 124 (shifts old rows 124-147 down to 125-148,
# which matches the committed diff exactly - each subsequent row's data
# slides down by one, and a brand-new entry appears at row 124).
$ws.Rows("124:124").Insert()

# Populate the new row 124 with a fresh weekly price entry (same as the
# former row 124 record, but for a later reporting date).
$ws.Range("A124").Value = 10
$ws.Range("B124").Value = "Vega Modelo de Temuco"
$ws.Range("C124").Value = "La Araucanía"
$ws.Range("D124").Value = 44476
$ws.Range("E124").Value = 9
$ws.Range("F124").Value = "Fruta"
$ws.Range("G124").Value = 100102
$ws.Range("H124").Value = "Cítricos"
$ws.Range("I124").Value = 100102006
$ws.Range("J124").Value = "Pomelo"
$ws.Range("K124").Value = "Start Ruby"
$ws.Range("L124").Value = "Primera"
$ws.Range("M124").Value = 90
$ws.Range("N124").Value = 12000
$ws.Range("O124").Value = 12000
$ws.Range("P124").Value = 12000
$ws.Range("Q124").Value = "$/bandeja 15 kilos granel"
$ws.Range("R124").Value = "Región de O'Higgins"
$ws.Range("S124").Value = 800
$ws.Range("T124").Value = 15
